# Add 15 new transaction rows (rows 65-79) for week ending 2021-01-24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Data layout: Row, Date(A), ReceiptNumber(B), Konto(C), Beskrivning(D), Debet(E), Kredit(F)
# Use $null to indicate "leave blank / do not set"
$rows = @(
    @(65, 44216, $null,    6570, "Pris banktjänster enligt faktura", 1250,   $null),
    @(66, 44216, $null,    $null, "Pris banktjänster enligt faktura", 0,      $null),
    @(67, 44216, $null,    1930, "Pris banktjänster enligt faktura", $null,  1250),
    @(68, 44217, 9211635, 3011, "Order 9211635 Swish +46707883566", $null,  1036.61),
    @(69, 44217, 9211635, 2611, "Order 9211635 Swish +46707883566", $null,  124.39),
    @(70, 44217, 9211635, 1930, "Order 9211635 Swish +46707883566", 1161,   $null),
    @(71, 44218, 5222159, 3011, "Order 5222159 Swish +46739582203", $null,  547.3200000000001),
    @(72, 44218, 5222159, 2611, "Order 5222159 Swish +46739582203", $null,  65.68000000000001),
    @(73, 44218, 5222159, 1930, "Order 5222159 Swish +46739582203", 613,    $null),
    @(74, 44218, $null,    4010, "M&S RB BROMMA K0135", 1529.97, $null),
    @(75, 44218, $null,    2645, "M&S RB BROMMA K0135", 183.59,  $null),
    @(76, 44218, $null,    1930, "M&S RB BROMMA K0135", $null,   1713.56),
    @(77, 44219, $null,    5670, "OKQ8 K0135", 764.26, $null),
    @(78, 44219, $null,    2641, "OKQ8 K0135", 191.07, $null),
    @(79, 44219, $null,    1930, "OKQ8 K0135", $null,  955.33)
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r[1]
    $cellA.NumberFormat = $dateFormat

    $cellB = $ws.Cells.Item($rowNum, 2)
    if ($null -eq $r[2]) { $cellB.Value = "" } else { $cellB.Value = $r[2] }

    $cellC = $ws.Cells.Item($rowNum, 3)
    if ($null -eq $r[3]) { $cellC.Value = "" } else { $cellC.Value = $r[3] }

    $ws.Cells.Item($rowNum, 4).Value = $r[4]

    $cellE = $ws.Cells.Item($rowNum, 5)
    if ($null -eq $r[5]) { $cellE.Value = "" } else { $cellE.Value = $r[5] }

    $cellF = $ws.Cells.Item($rowNum, 6)
    if ($null -eq $r[6]) { $cellF.Value = "" } else { $cellF.Value = $r[6] }
}
